$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 91.25
$ws.Range("I2").Value = 91.25
$ws.Range("K2").Value = 91.25
$ws.Range("M2").Value = 21.75

$ws.Range("H40").Value = 114908.96
$ws.Range("J40").Value = 3956.75
$ws.Range("L40").Value = 3956.75
$ws.Range("N40").Value = -4306.75

$ws.Range("H53").Value = 1258.6072
$ws.Range("I53").Value = 1186.7333
$ws.Range("K53").Value = 1186.7333
$ws.Range("M53").Value = -549.7333000000001

$ws.Range("H86").Value = 181824900
$ws.Range("I86").Value = 142864420
$ws.Range("J86").Value = 250005740
$ws.Range("K86").Value = 142864420
$ws.Range("L86").Value = 250005740
$ws.Range("M86").Value = -142863297
$ws.Range("N86").Value = -250007986

$ws.Range("H89").Value = 181824900
$ws.Range("I89").Value = 142864420
$ws.Range("J89").Value = 250005740
$ws.Range("K89").Value = 714322100
$ws.Range("L89").Value = 1250028700
$ws.Range("M89").Value = -714316484
$ws.Range("N89").Value = -1250039932

$ws.Range("H113").Value = 1953.5416
$ws.Range("J113").Value = 2700.8572
$ws.Range("L113").Value = 2700.8572
$ws.Range("N113").Value = -9208.8572

$ws.Range("H114").Value = 91316.336
$ws.Range("J114").Value = 91316.336
$ws.Range("L114").Value = 91316.336
$ws.Range("N114").Value = -99994.336

$ws.Range("H131").Value = 5006.9
$ws.Range("I131").Value = 1033.8
$ws.Range("K131").Value = 3101.4
$ws.Range("M131").Value = 1938.6

$ws.Range("H132").Value = 3667.5918
$ws.Range("I132").Value = 4230.7427
$ws.Range("K132").Value = 12692.2281
$ws.Range("M132").Value = -10162.2281

$ws.Range("H141").Value = 3323.5833
$ws.Range("I141").Value = 2408.6
$ws.Range("K141").Value = 7225.799999999999
$ws.Range("M141").Value = -2045.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 664.7895
$ws.Range("I2").Value = 514.5625
$ws.Range("K2").Value = 514.5625
$ws.Range("M2").Value = -401.5625

$ws.Range("H11").Value = 333996.66
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 333996.66
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 333996.66
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -334284.66

$ws.Range("H61").Value = 3016.3408
$ws.Range("I61").Value = 2843.2163
$ws.Range("J61").Value = 3931.4285
$ws.Range("K61").Value = 2843.2163
$ws.Range("L61").Value = 3931.4285
$ws.Range("M61").Value = -2631.2163
$ws.Range("N61").Value = -4355.4285

$ws.Range("H63").Value = 143960.4
$ws.Range("I63").Value = 3001.6667
$ws.Range("K63").Value = 3001.6667
$ws.Range("M63").Value = -2315.6667

$ws.Range("H66").Value = 143960.4
$ws.Range("I66").Value = 3001.6667
$ws.Range("K66").Value = 15008.3335
$ws.Range("M66").Value = -11576.3335

$ws.Range("H102").Value = 1778
$ws.Range("I102").Value = 1778
$ws.Range("K102").Value = 1778
$ws.Range("M102").Value = -156

$ws.Range("H105").Value = 119994
$ws.Range("J105").Value = 119994
$ws.Range("L105").Value = 119994
$ws.Range("N105").Value = -126982

$ws.Range("H107").Value = 37984.5
$ws.Range("J107").Value = 37984.5
$ws.Range("L107").Value = 37984.5
$ws.Range("N107").Value = -45664.5

$ws.Range("H111").Value = 79999
$ws.Range("I111").Value = 79999
$ws.Range("J111").Value = 79999
$ws.Range("K111").Value = 79999
$ws.Range("L111").Value = 79999
$ws.Range("M111").Value = -75909
$ws.Range("N111").Value = -88179

$ws.Range("H113").Value = 48994
$ws.Range("J113").Value = 48994
$ws.Range("L113").Value = 48994
$ws.Range("N113").Value = -57672

$ws.Range("H115").Value = 50996.168
$ws.Range("I115").Value = 33333
$ws.Range("J115").Value = 68659.336
$ws.Range("K115").Value = 33333
$ws.Range("L115").Value = 68659.336
$ws.Range("M115").Value = -31766
$ws.Range("N115").Value = -71793.336

$ws.Range("H116").Value = 664.7895
$ws.Range("I116").Value = 514.5625
$ws.Range("K116").Value = 514.5625
$ws.Range("M116").Value = 1779.4375

$ws.Range("H131").Value = 64165.668
$ws.Range("J131").Value = 42499
$ws.Range("L131").Value = 42499
$ws.Range("N131").Value = -52579

$ws.Range("H132").Value = 3683.1667
$ws.Range("I132").Value = 3209
$ws.Range("K132").Value = 9627
$ws.Range("M132").Value = -7097

$ws.Range("H136").Value = 3016.3408
$ws.Range("I136").Value = 2843.2163
$ws.Range("J136").Value = 3931.4285
$ws.Range("K136").Value = 8529.6489
$ws.Range("L136").Value = 11794.2855
$ws.Range("M136").Value = -5979.6489
$ws.Range("N136").Value = -16894.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 664.7895
$ws.Range("I3").Value = 514.5625
$ws.Range("K3").Value = 514.5625
$ws.Range("M3").Value = -400.5625

$ws.Range("H20").Value = 32688.414
$ws.Range("I20").Value = 44334.523
$ws.Range("J20").Value = 2117.375
$ws.Range("K20").Value = 44334.523
$ws.Range("L20").Value = 2117.375
$ws.Range("M20").Value = -44087.523
$ws.Range("N20").Value = -2611.375

$ws.Range("H50").Value = 51375.75
$ws.Range("J50").Value = 40931.332
$ws.Range("L50").Value = 40931.332
$ws.Range("N50").Value = -42079.332

$ws.Range("H55").Value = 64992
$ws.Range("J55").Value = 64992
$ws.Range("L55").Value = 64992
$ws.Range("N55").Value = -65538

$ws.Range("H80").Value = 694.5
$ws.Range("J80").Value = 679.6667
$ws.Range("L80").Value = 679.6667
$ws.Range("N80").Value = -2675.6667

$ws.Range("H83").Value = 694.5
$ws.Range("J83").Value = 679.6667
$ws.Range("L83").Value = 3398.3335
$ws.Range("N83").Value = -13382.3335

$ws.Range("H86").Value = 3637.5173
$ws.Range("I86").Value = 2499.476
$ws.Range("K86").Value = 2499.476
$ws.Range("M86").Value = -1376.476

$ws.Range("H89").Value = 3637.5173
$ws.Range("I89").Value = 2499.476
$ws.Range("K89").Value = 12497.38
$ws.Range("M89").Value = -6881.380000000001

$ws.Range("H94").Value = 1200.2727
$ws.Range("I94").Value = 301.44446
$ws.Range("J94").Value = 5245
$ws.Range("K94").Value = 301.44446
$ws.Range("L94").Value = 5245
$ws.Range("M94").Value = 149.55554
$ws.Range("N94").Value = -6147

$ws.Range("H111").Value = 79990
$ws.Range("J111").Value = 79990
$ws.Range("L111").Value = 79990
$ws.Range("N111").Value = -88170

$ws.Range("H112").Value = 79990
$ws.Range("J112").Value = 79990
$ws.Range("L112").Value = 79990
$ws.Range("N112").Value = -82944

$ws.Range("H114").Value = 134963
$ws.Range("J114").Value = 134963
$ws.Range("L114").Value = 134963
$ws.Range("N114").Value = -143641

$ws.Range("H134").Value = 20839100
$ws.Range("I134").Value = 3035.1667
$ws.Range("J134").Value = 33340738
$ws.Range("K134").Value = 9105.500100000001
$ws.Range("L134").Value = 100022214
$ws.Range("M134").Value = -6570.500100000001
$ws.Range("N134").Value = -100027284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1833.24
$ws.Range("I16").Value = 1692.3182
$ws.Range("K16").Value = 1692.3182
$ws.Range("M16").Value = -1405.3182

$ws.Range("H31").Value = 6144.5
$ws.Range("I31").Value = 2085.1875
$ws.Range("J31").Value = 10783.714
$ws.Range("K31").Value = 2085.1875
$ws.Range("L31").Value = 10783.714
$ws.Range("M31").Value = -1790.1875
$ws.Range("N31").Value = -11373.714

$ws.Range("H34").Value = 6144.5
$ws.Range("I34").Value = 2085.1875
$ws.Range("J34").Value = 10783.714
$ws.Range("K34").Value = 2085.1875
$ws.Range("L34").Value = 10783.714
$ws.Range("M34").Value = -1883.1875
$ws.Range("N34").Value = -11187.714

$ws.Range("H58").Value = 3346.2666
$ws.Range("J58").Value = 3399.8462
$ws.Range("L58").Value = 3399.8462
$ws.Range("N58").Value = -3805.8462

$ws.Range("H62").Value = 2409.75
$ws.Range("I62").Value = 2409.75
$ws.Range("K62").Value = 2409.75
$ws.Range("M62").Value = -1785.75

$ws.Range("H65").Value = 2409.75
$ws.Range("I65").Value = 2409.75
$ws.Range("K65").Value = 12048.75
$ws.Range("M65").Value = -8928.75

$ws.Range("H99").Value = 2088.375
$ws.Range("I99").Value = 1053
$ws.Range("J99").Value = 3123.75
$ws.Range("K99").Value = 1053
$ws.Range("L99").Value = 3123.75
$ws.Range("M99").Value = 445
$ws.Range("N99").Value = -6119.75

$ws.Range("H100").Value = 113850
$ws.Range("J100").Value = 113850
$ws.Range("L100").Value = 113850
$ws.Range("N100").Value = -116014

$ws.Range("H113").Value = 1833.24
$ws.Range("I113").Value = 1692.3182
$ws.Range("K113").Value = 1692.3182
$ws.Range("M113").Value = 477.6818000000001

$ws.Range("H126").Value = 2088.375
$ws.Range("I126").Value = 1053
$ws.Range("J126").Value = 3123.75
$ws.Range("K126").Value = 3159
$ws.Range("L126").Value = 9371.25
$ws.Range("M126").Value = -689
$ws.Range("N126").Value = -14311.25

$ws.Range("H132").Value = 2224.8462
$ws.Range("I132").Value = 1510.2
$ws.Range("K132").Value = 4530.6
$ws.Range("M132").Value = -2000.6

$ws.Range("H136").Value = 3346.2666
$ws.Range("J136").Value = 3399.8462
$ws.Range("L136").Value = 10199.5386
$ws.Range("N136").Value = -15299.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 879.05884
$ws.Range("J5").Value = 1220.125
$ws.Range("L5").Value = 3660.375
$ws.Range("N5").Value = -3884.375

$ws.Range("H51").Value = 3707.3333
$ws.Range("J51").Value = 4360
$ws.Range("L51").Value = 13080
$ws.Range("N51").Value = -14000

$ws.Range("H132").Value = 1825.9166
$ws.Range("I132").Value = 1001.8571
$ws.Range("J132").Value = 2979.6
$ws.Range("K132").Value = 9016.713899999999
$ws.Range("L132").Value = 26816.4
$ws.Range("M132").Value = -6486.713899999999
$ws.Range("N132").Value = -31876.4

$ws.Range("H134").Value = 1306.6923
$ws.Range("I134").Value = 1306.6923
$ws.Range("K134").Value = 3920.0769
$ws.Range("M134").Value = 1149.9231

$ws.Range("H135").Value = 879.05884
$ws.Range("J135").Value = 1220.125
$ws.Range("L135").Value = 10981.125
$ws.Range("N135").Value = -16051.125

$ws.Range("H137").Value = 1567.375
$ws.Range("I137").Value = 1567.375
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4702.125
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 397.875
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 37875
$ws.Range("J20").Value = 17666
$ws.Range("L20").Value = 17666
$ws.Range("N20").Value = -18156

$ws.Range("H24").Value = 5755713
$ws.Range("J24").Value = 57998.4
$ws.Range("L24").Value = 57998.4
$ws.Range("N24").Value = -58344.4

$ws.Range("H70").Value = 21858.354
$ws.Range("I70").Value = 63599.4
$ws.Range("J70").Value = 4466.25
$ws.Range("K70").Value = 63599.4
$ws.Range("L70").Value = 4466.25
$ws.Range("M70").Value = -63329.4
$ws.Range("N70").Value = -5006.25

$ws.Range("H73").Value = 21858.354
$ws.Range("I73").Value = 63599.4
$ws.Range("J73").Value = 4466.25
$ws.Range("K73").Value = 63599.4
$ws.Range("L73").Value = 4466.25
$ws.Range("M73").Value = -62663.4
$ws.Range("N73").Value = -6338.25

$ws.Range("H97").Value = 1688.1111
$ws.Range("I97").Value = 1099.1428
$ws.Range("K97").Value = 1099.1428
$ws.Range("M97").Value = -603.1428000000001

$ws.Range("H102").Value = 2427.1428
$ws.Range("I102").Value = 2427.1428
$ws.Range("K102").Value = 2427.1428
$ws.Range("M102").Value = -805.1428000000001

$ws.Range("H114").Value = 95656.336
$ws.Range("I114").Value = 153000
$ws.Range("J114").Value = 66984.5
$ws.Range("K114").Value = 153000
$ws.Range("L114").Value = 66984.5
$ws.Range("M114").Value = -148661
$ws.Range("N114").Value = -75662.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3117.25
$ws.Range("I7").Value = 3117.25
$ws.Range("K7").Value = 3117.25
$ws.Range("M7").Value = -3005.25

$ws.Range("H46").Value = 2361.5789
$ws.Range("I46").Value = 532.4286
$ws.Range("K46").Value = 532.4286
$ws.Range("M46").Value = -344.4286

$ws.Range("H59").Value = 85587.86
$ws.Range("J59").Value = 85587.86
$ws.Range("L59").Value = 85587.86
$ws.Range("N59").Value = -86895.86

$ws.Range("H61").Value = 1810.591
$ws.Range("I61").Value = 1677.3
$ws.Range("K61").Value = 1677.3
$ws.Range("M61").Value = -1475.3

$ws.Range("H68").Value = 2508.2307
$ws.Range("I68").Value = 2405
$ws.Range("J68").Value = 2628.6667
$ws.Range("K68").Value = 2405
$ws.Range("L68").Value = 2628.6667
$ws.Range("M68").Value = -1656
$ws.Range("N68").Value = -4126.6667

$ws.Range("H71").Value = 2508.2307
$ws.Range("I71").Value = 2405
$ws.Range("J71").Value = 2628.6667
$ws.Range("K71").Value = 12025
$ws.Range("L71").Value = 13143.3335
$ws.Range("M71").Value = -8281
$ws.Range("N71").Value = -20631.3335

$ws.Range("H105").Value = 75500
$ws.Range("J105").Value = 75500
$ws.Range("L105").Value = 75500
$ws.Range("N105").Value = -82488

$ws.Range("H113").Value = 1810.591
$ws.Range("I113").Value = 1677.3
$ws.Range("K113").Value = 1677.3
$ws.Range("M113").Value = 492.7

$ws.Range("H114").Value = 49429.668
$ws.Range("J114").Value = 49429.668
$ws.Range("L114").Value = 49429.668
$ws.Range("N114").Value = -58107.668

$ws.Range("H115").Value = 80783.8
$ws.Range("J115").Value = 80783.8
$ws.Range("L115").Value = 80783.8
$ws.Range("N115").Value = -83133.8

$ws.Range("H120").Value = 74296.336
$ws.Range("J120").Value = 74296.336
$ws.Range("L120").Value = 74296.336
$ws.Range("N120").Value = -83972.336

$ws.Range("H122").Value = 4263.75
$ws.Range("I122").Value = 4016.6667
$ws.Range("K122").Value = 12050.0001
$ws.Range("M122").Value = -9600.000100000001

$ws.Range("H126").Value = 3117.25
$ws.Range("I126").Value = 3117.25
$ws.Range("K126").Value = 9351.75
$ws.Range("M126").Value = -6881.75

$ws.Range("H128").Value = 89163.164
$ws.Range("J128").Value = 89163.164
$ws.Range("L128").Value = 89163.164
$ws.Range("N128").Value = -99123.164

$ws.Range("H136").Value = 4571.4287
$ws.Range("J136").Value = 5166.6665
$ws.Range("L136").Value = 15499.9995
$ws.Range("N136").Value = -20599.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 50544.668
$ws.Range("J27").Value = 50544.668
$ws.Range("L27").Value = 50544.668
$ws.Range("N27").Value = -50682.668

$ws.Range("H31").Value = 59019
$ws.Range("J31").Value = 59019
$ws.Range("L31").Value = 59019
$ws.Range("N31").Value = -59715

$ws.Range("H62").Value = 6186.9165
$ws.Range("I62").Value = 5932.75
$ws.Range("J62").Value = 6314
$ws.Range("K62").Value = 5932.75
$ws.Range("L62").Value = 6314
$ws.Range("M62").Value = -5308.75
$ws.Range("N62").Value = -7562

$ws.Range("H65").Value = 6186.9165
$ws.Range("I65").Value = 5932.75
$ws.Range("J65").Value = 6314
$ws.Range("K65").Value = 29663.75
$ws.Range("L65").Value = 31570
$ws.Range("M65").Value = -26543.75
$ws.Range("N65").Value = -37810

$ws.Range("H81").Value = 3963.611
$ws.Range("J81").Value = 4192.769
$ws.Range("L81").Value = 8385.538
$ws.Range("N81").Value = -10507.538

$ws.Range("H84").Value = 3963.611
$ws.Range("J84").Value = 4192.769
$ws.Range("L84").Value = 41927.69
$ws.Range("N84").Value = -52535.69

$ws.Range("H110").Value = 60446.375
$ws.Range("J110").Value = 60446.375
$ws.Range("L110").Value = 60446.375
$ws.Range("N110").Value = -68626.375

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

$ws.Range("H118").Value = 109499.5
$ws.Range("J118").Value = 109499.5
$ws.Range("L118").Value = 109499.5
$ws.Range("N118").Value = -112813.5

$ws.Range("H132").Value = 4949.7393
$ws.Range("I132").Value = 4562.9414
$ws.Range("K132").Value = 13688.8242
$ws.Range("M132").Value = -11158.8242
